# "Full Testing.xlsx" - add the new "search box" test cases to the bottom
# of the testing sheet (rows 71-76), continuing the existing numbered test
# list (tests 49-54) with the same box-border / wrap-text formatting used
# by the rows above them.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 71 already exists (blank, but pre-formatted) - just fill in its values.
$ws.Range("A71").Value = 49
$ws.Range("B71").Value = "is the search box working"

# Rows 72-76 are brand new. Copy row 71's formatting (border box + wrap
# text + font) down into them before filling in their values, so they
# match the look of the rest of the test list.
$ws.Range("A71:H71").Copy($ws.Range("A72:H76"))

$ws.Range("A72").Value = 50
$ws.Range("B72").Value = "Does the search box prevent XSS attacks"

$ws.Range("A73").Value = 51
$ws.Range("B73").Value = "Does the search box prevent SQL injections?"

$ws.Range("A74").Value = 52
$ws.Range("B74").Value = "is the search box working"

$ws.Range("A75").Value = 53
$ws.Range("B75").Value = "is the search box working"

$ws.Range("A76").Value = 54
$ws.Range("B76").Value = "is the search box working"

# Row heights roughly matching the wrapped text in each row.
$ws.Rows.Item(72).RowHeight = 36
$ws.Rows.Item(73).RowHeight = 25.5
$ws.Rows.Item(74).RowHeight = 15.75
$ws.Rows.Item(75).RowHeight = 15.75
$ws.Rows.Item(76).RowHeight = 15.75

# Restore the view state left behind in the workbook (active sheet,
# selection, scroll position) from the last edit session.
$ws.Activate()
$ws.Range("A29:XFD29").Select()
